$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.405.46"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.501.49"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.16"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.47"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  +6.08%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("E11").Value = "  +4.24%  "
$ws.Range("D12").Value = "4.097.71"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D15").Value = "3.508.41"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.76"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").Value = "64.397.86"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.61"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.583"
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("D23").Value = "3.642.22"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.48"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.67"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.43"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.20"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.156"
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("D34").Value = "3.529.93"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.33"
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "165.51"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.44"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.02"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.18"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.929"
$ws.Range("E48").Value = "  +4.11%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.422.60"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0259"
$ws.Range("E51").Value = "  +0.56%  "
